$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add a new entry (LeetCode 1283 - "Find the Smallest Divisor Given a
# Threshold") as a new row 19. The previously-blank spacer row that used to
# occupy row 19 is replaced by this real data row; everything below keeps
# its existing row number.
# ---------------------------------------------------------------------------

# Start from row 18's formatting (same "visual group" of rows - header style
# s=32, fills, borders, wrap text, etc.) and apply it to row 19 first so the
# new cells inherit the correct look before we fill in values.
$ws.Range("A18:F18").Copy()
$ws.Range("A19:F19").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the new row's data. Question text (C19) is set before the
# category text (B19) so the new shared strings are appended in the same
# order as the target workbook (Question text first, then category).
$ws.Range("A19").Value = 1283
$ws.Range("C19").Value = "Find the Smallest Divisor Given a Threshold"
$ws.Range("B19").Value = "CN/LC/GFG"
$ws.Range("D19").Value = "Java"
$ws.Range("E19").Value = "Medium"

# Add the LeetCode hyperlink on the question cell, matching the existing
# hyperlinked question cells elsewhere in the sheet.
$ws.Hyperlinks.Add($ws.Range("C19"), "https://leetcode.com/problems/find-the-smallest-divisor-given-a-threshold/")

# Re-apply the row 18 cell formatting once more so the hyperlink-insertion
# side effect (which forces its own "hyperlink" look onto the cell) is
# overwritten back to the normal style used by every other question cell.
$ws.Range("A18:F18").Copy()
$ws.Range("A19:F19").PasteSpecial(-4122)  # xlPasteFormats

# Re-set the values (PasteSpecial(xlPasteFormats) only touches formatting,
# but do this defensively in case paste order ever changes) and keep F19
# empty, matching the target.
$ws.Range("A19").Value = 1283
$ws.Range("D19").Value = "Java"
$ws.Range("E19").Value = "Medium"

# ---------------------------------------------------------------------------
# Column E ("Level") fill-colour fix: E17 and E18 ("Medium") were using a
# mismatched fill colour compared to every other "Medium" cell in the
# column (e.g. E7:E11, E13, E14). Copy the correct format from one of those
# cells onto E17, E18 and the newly-added E19.
# ---------------------------------------------------------------------------
$ws.Range("E7").Copy()
$ws.Range("E17").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E18").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E19").PasteSpecial(-4122)  # xlPasteFormats

# Re-assert the text values for E17/E18/E19 since PasteSpecial(xlPasteFormats)
# only changes formatting, not content (values already correct, kept for
# clarity/safety).
$ws.Range("E17").Value = "Medium"
$ws.Range("E18").Value = "Medium"
$ws.Range("E19").Value = "Medium"

# ---------------------------------------------------------------------------
# Update the sheet's active selection to match the post-edit state recorded
# in the workbook (user ended up with C23 selected after inserting the new
# row).
# ---------------------------------------------------------------------------
$ws.Range("C23").Select()
